# Generate Report for Handback
#
# Two e2e test markdown files were regenerated with new GUIDs/timestamps:
#   7e78e3b5-ac59-4e94-8ba5-8a9a2d0cd8b6  ->  1ee4c6c2-6e30-4c61-980c-0330ef322f42
#   843ba29b-8ddd-4251-a1d3-4da9f64c97bc  ->  ffff36796dbf-10b5-45c6-9416-82b89de0df2a
# and the corresponding handoff/handback xliff file names + timestamps were
# refreshed as well. This script updates every cell (and hyperlink display
# text) across the three report sheets (Overview, zh-cn, de-de) to reflect
# the new handback run, while leaving the hyperlink target URLs untouched
# (they still point at the original commit).

$wb = $excel.ActiveWorkbook

$oldGuid1 = "7e78e3b5-ac59-4e94-8ba5-8a9a2d0cd8b6"
$newGuid1 = "1ee4c6c2-6e30-4c61-980c-0330ef322f42"
$oldGuid2 = "843ba29b-8ddd-4251-a1d3-4da9f64c97bc"
$newGuid2 = "ffff36796dbf-10b5-45c6-9416-82b89de0df2a"

$oldHash1 = "0ea02e6897cec90f1ccef63300de1c79b858b989"
$newHash1 = "d775bd04f86a5c9438bee9800cc07797defe2932"

# -----------------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-09-01 23:08:15"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-09-01 23:08:15"

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Row -eq 2) {
        $h.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($h.Range.Row -eq 3) {
        $h.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# -----------------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 23:08:09"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-01 23:08:36"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 23:08:09"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-01 23:08:36"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Row -eq 2) {
        $h.TextToDisplay = "$newGuid1.md"
    } elseif ($h.Range.Row -eq 3) {
        $h.TextToDisplay = "$newGuid2.md"
    }
}

# -----------------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 23:08:15"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-01 23:08:44"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 23:08:15"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-01 23:08:44"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Row -eq 2) {
        $h.TextToDisplay = "$newGuid1.md"
    } elseif ($h.Range.Row -eq 3) {
        $h.TextToDisplay = "$newGuid2.md"
    }
}
